$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in rows 29 and 30 with the new "World and Machine" / "Goals" entries
$ws.Range("A29").Value = 43762
$ws.Range("B29").Value = "World and Machine"
$ws.Range("C29").Value = 1

$ws.Range("A30").Value = 40841
$ws.Range("B30").Value = "Goals"
$ws.Range("C30").Value = 1

# Update the view so the sheet opens scrolled to the new rows, matching the
# author's last on-screen selection when the edit was made.
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("G27").Select()
